# Edit presentation per the target change:
#  1) Slide 1, 3rd shape ("Content Placeholder 4" textbox with the team
#     member list): remove the name "Sirul" (and the following space)
#     from the paragraph "Vaishnavi Sirul Velaga ", leaving
#     "Vaishnavi Velaga ".
#  2) Slide 7, 2nd shape ("Content Placeholder 2"): add a new "Hints"
#     paragraph right after the existing "Role of GUI" paragraph, before
#     the trailing empty paragraph.

$p = $ppt.ActivePresentation

# --- 1) Slide 1: "Vaishnavi Sirul Velaga " -> "Vaishnavi Velaga " ---
$s1 = $p.Slides.Item(1)
$shape1 = $s1.Shapes.Item(3)
$tr1 = $shape1.TextFrame.TextRange

# Locate the "Sirul " substring (including the trailing space that
# separated it from "Velaga") and delete it.
$sirul = $tr1.Characters(24, 6)
$sirul.Text = ""

# --- 2) Slide 7: insert a new "Hints" paragraph after "Role of GUI" ---
$s7 = $p.Slides.Item(7)
$shape7 = $s7.Shapes.Item(2)
$tr7 = $shape7.TextFrame.TextRange

# Appends "Hints" + a paragraph break right after the existing text,
# i.e. right after "Role of GUI" and before the already-present
# trailing empty paragraph.
[void]$tr7.InsertAfter("Hints`r")

# Match the font size used by the rest of the bullet list (28pt).
$newPara = $tr7.Characters(26, 5)
$newPara.Font.Size = 28
